$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2193
$ws.Range("I62").Value = 2289.6
$ws.Range("J62").Value = 1951.5
$ws.Range("K62").Value = 2289.6
$ws.Range("L62").Value = 1951.5
$ws.Range("M62").Value = -1665.6
$ws.Range("N62").Value = -3199.5

$ws.Range("H65").Value = 2193
$ws.Range("I65").Value = 2289.6
$ws.Range("J65").Value = 1951.5
$ws.Range("K65").Value = 11448
$ws.Range("L65").Value = 9757.5
$ws.Range("M65").Value = -8328
$ws.Range("N65").Value = -15997.5

$ws.Range("H98").Value = 2045.4166
$ws.Range("I98").Value = 2131.0527
$ws.Range("J98").Value = 1720
$ws.Range("K98").Value = 2131.0527
$ws.Range("L98").Value = 1720
$ws.Range("M98").Value = -633.0527000000002
$ws.Range("N98").Value = -4716

$ws.Range("H122").Value = 2045.4166
$ws.Range("I122").Value = 2131.0527
$ws.Range("J122").Value = 1720
$ws.Range("K122").Value = 6393.158100000001
$ws.Range("L122").Value = 5160
$ws.Range("M122").Value = -3943.158100000001
$ws.Range("N122").Value = -10060

$ws.Range("H129").Value = 91968.91
$ws.Range("J129").Value = 168166.33
$ws.Range("L129").Value = 504498.99
$ws.Range("N129").Value = -514498.99

$ws.Range("H137").Value = 2237.8462
$ws.Range("I137").Value = 2214.1428
$ws.Range("K137").Value = 6642.428400000001
$ws.Range("M137").Value = -4092.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 109822.88
$ws.Range("J63").Value = 1000
$ws.Range("L63").Value = 1000
$ws.Range("N63").Value = -2372

$ws.Range("H66").Value = 109822.88
$ws.Range("J66").Value = 1000
$ws.Range("L66").Value = 5000
$ws.Range("N66").Value = -11864

$ws.Range("H80").Value = 27200
$ws.Range("J80").Value = 27200
$ws.Range("L80").Value = 27200
$ws.Range("N80").Value = -29196

$ws.Range("H83").Value = 27200
$ws.Range("J83").Value = 27200
$ws.Range("L83").Value = 81600
$ws.Range("N83").Value = -91584

$ws.Range("H132").Value = 4873.357
$ws.Range("I132").Value = 4852.25
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 14556.75
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -12026.75
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 31776.857
$ws.Range("I82").Value = 3257
$ws.Range("J82").Value = 36530.168
$ws.Range("K82").Value = 3257
$ws.Range("L82").Value = 36530.168
$ws.Range("M82").Value = -2874
$ws.Range("N82").Value = -37296.168

$ws.Range("H85").Value = 31776.857
$ws.Range("I85").Value = 3257
$ws.Range("J85").Value = 36530.168
$ws.Range("K85").Value = 3257
$ws.Range("L85").Value = 36530.168
$ws.Range("M85").Value = -1931
$ws.Range("N85").Value = -39182.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 406.36365
$ws.Range("I22").Value = 397
$ws.Range("K22").Value = 397
$ws.Range("M22").Value = -47

$ws.Range("H31").Value = 8002816.5
$ws.Range("I31").Value = 2887.6956
$ws.Range("K31").Value = 2887.6956
$ws.Range("M31").Value = -2592.6956

$ws.Range("H34").Value = 8002816.5
$ws.Range("I34").Value = 2887.6956
$ws.Range("K34").Value = 2887.6956
$ws.Range("M34").Value = -2685.6956

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 562.2692
$ws.Range("I5").Value = 422.78262
$ws.Range("K5").Value = 1268.34786
$ws.Range("M5").Value = -1156.34786

$ws.Range("H122").Value = 1617.6522
$ws.Range("J122").Value = 1386.8667
$ws.Range("L122").Value = 12481.8003
$ws.Range("N122").Value = -17381.8003

$ws.Range("H131").Value = 2004691.9
$ws.Range("J131").Value = 2670170
$ws.Range("L131").Value = 8010510
$ws.Range("N131").Value = -8020590

$ws.Range("H132").Value = 1339.4
$ws.Range("J132").Value = 2407.7273
$ws.Range("L132").Value = 21669.5457
$ws.Range("N132").Value = -26729.5457

$ws.Range("H135").Value = 562.2692
$ws.Range("I135").Value = 422.78262
$ws.Range("K135").Value = 3805.04358
$ws.Range("M135").Value = -1270.04358

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1630.1428
$ws.Range("I102").Value = 1532.1538
$ws.Range("J102").Value = 1789.375
$ws.Range("K102").Value = 1532.1538
$ws.Range("L102").Value = 1789.375
$ws.Range("M102").Value = 89.84619999999995
$ws.Range("N102").Value = -5033.375

$ws.Range("H113").Value = 22728326
$ws.Range("J113").Value = 1192
$ws.Range("L113").Value = 1192
$ws.Range("N113").Value = -5532

$ws.Range("H122").Value = 1541.3478
$ws.Range("I122").Value = 540.2727
$ws.Range("J122").Value = 2459
$ws.Range("K122").Value = 1620.8181
$ws.Range("L122").Value = 7377
$ws.Range("M122").Value = 829.1819
$ws.Range("N122").Value = -12277

$ws.Range("H126").Value = 6601.2
$ws.Range("I126").Value = 7002.4
$ws.Range("J126").Value = 6200
$ws.Range("K126").Value = 21007.2
$ws.Range("L126").Value = 18600
$ws.Range("M126").Value = -18537.2
$ws.Range("N126").Value = -23540

$ws.Range("H132").Value = 89962.35000000001
$ws.Range("I132").Value = 202433.6
$ws.Range("J132").Value = 3446
$ws.Range("K132").Value = 607300.8
$ws.Range("L132").Value = 10338
$ws.Range("M132").Value = -604770.8
$ws.Range("N132").Value = -15398

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2389
$ws.Range("I40").Value = 1980.2858
$ws.Range("J40").Value = 5250
$ws.Range("K40").Value = 1980.2858
$ws.Range("L40").Value = 5250
$ws.Range("M40").Value = -1844.2858
$ws.Range("N40").Value = -5522

$ws.Range("H47").Value = 17866.666
$ws.Range("J47").Value = 17866.666
$ws.Range("L47").Value = 17866.666
$ws.Range("N47").Value = -18846.666

$ws.Range("H52").Value = 17866.666
$ws.Range("J52").Value = 17866.666
$ws.Range("L52").Value = 17866.666
$ws.Range("N52").Value = -18332.666

$ws.Range("H122").Value = 1815.3846
$ws.Range("I122").Value = 1666.6666
$ws.Range("J122").Value = 1860
$ws.Range("K122").Value = 4999.9998
$ws.Range("L122").Value = 5580
$ws.Range("M122").Value = -2549.9998
$ws.Range("N122").Value = -10480
